# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period detail table (rows 16-36, columns C:G) is re-sorted:
# previously grouped by worker (outer) then by period descending (inner),
# now grouped by period ascending (outer) then by worker (inner). The
# "Valor Mora" (F) amount of 58667 now always accompanies period 2405 for
# every worker (it used to be tied to whichever row came first per
# worker). The "Salario Basico" (G) for JOSE GUILLERMO ANGULO VIAÑA is
# updated from 1800000 to 2000000 for every one of his rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# doc type, doc number, worker name, period, valor mora, salario basico
$data = @(
    @("CC", "45546839",   "CLAUDIA PATRICIA MEJIA RAMIREZ", "2311", 80000, 2000000),
    @("CC", "1047422454", "MARTICELA BAYTER DORIA",         "2311", 80000, 2000000),
    @("CC", "1143360875", "JOSE GUILLERMO ANGULO VIAÑA",    "2311", 80000, 2000000),

    @("CC", "45546839",   "CLAUDIA PATRICIA MEJIA RAMIREZ", "2312", 80000, 2000000),
    @("CC", "1047422454", "MARTICELA BAYTER DORIA",         "2312", 80000, 2000000),
    @("CC", "1143360875", "JOSE GUILLERMO ANGULO VIAÑA",    "2312", 80000, 2000000),

    @("CC", "45546839",   "CLAUDIA PATRICIA MEJIA RAMIREZ", "2401", 80000, 2000000),
    @("CC", "1047422454", "MARTICELA BAYTER DORIA",         "2401", 80000, 2000000),
    @("CC", "1143360875", "JOSE GUILLERMO ANGULO VIAÑA",    "2401", 80000, 2000000),

    @("CC", "45546839",   "CLAUDIA PATRICIA MEJIA RAMIREZ", "2402", 80000, 2000000),
    @("CC", "1047422454", "MARTICELA BAYTER DORIA",         "2402", 80000, 2000000),
    @("CC", "1143360875", "JOSE GUILLERMO ANGULO VIAÑA",    "2402", 80000, 2000000),

    @("CC", "45546839",   "CLAUDIA PATRICIA MEJIA RAMIREZ", "2403", 80000, 2000000),
    @("CC", "1047422454", "MARTICELA BAYTER DORIA",         "2403", 80000, 2000000),
    @("CC", "1143360875", "JOSE GUILLERMO ANGULO VIAÑA",    "2403", 80000, 2000000),

    @("CC", "45546839",   "CLAUDIA PATRICIA MEJIA RAMIREZ", "2404", 80000, 2000000),
    @("CC", "1047422454", "MARTICELA BAYTER DORIA",         "2404", 80000, 2000000),
    @("CC", "1143360875", "JOSE GUILLERMO ANGULO VIAÑA",    "2404", 80000, 2000000),

    @("CC", "45546839",   "CLAUDIA PATRICIA MEJIA RAMIREZ", "2405", 58667, 2000000),
    @("CC", "1047422454", "MARTICELA BAYTER DORIA",         "2405", 58667, 2000000),
    @("CC", "1143360875", "JOSE GUILLERMO ANGULO VIAÑA",    "2405", 58667, 2000000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 2).Value = $rec[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $rec[1]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $rec[2]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $rec[3]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $rec[4]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value = $rec[5]   # G - Salario Basico
}
